# Excel COM-interop edit script
# - Swap columns A ("brand") and B ("name") so that A holds name values and
#   B holds brand values (header + the four existing data rows).
# - Add a new row 5: same brand/product/order/qty/status as row 4, but with
#   a new name value "테스트이름".
# - The custom column-width formatting that used to target column B
#   (name) now targets column A (name's new home).
# - Move the active-cell selection to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): A<->B swap: brand/name -> name/brand ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "brand"

# --- Row 2 ---
$ws.Range("A2").Value = "홍길동"
$ws.Range("B2").Value = "보니토"

# --- Row 3 ---
$ws.Range("A3").Value = "김철수"
$ws.Range("B3").Value = "보니토"

# --- Row 4 ---
$ws.Range("A4").Value = "루나"
$ws.Range("B4").Value = "보니토"

# --- New row 5 (copy of row 4's brand/product/date/qty/status, new name) ---
$ws.Range("A5").Value = "테스트이름"
$ws.Range("B5").Value = "보니토"
$ws.Range("C5").Value = "이거"
$ws.Range("D5").Value = "저거"
$ws.Range("E5").Value = "시러"
$ws.Range("F5").Value = "ㄹㄴㅇ"

# C4 carries a wrap-text style; mirror it on the new C5 cell
$ws.Range("C5").WrapText = $true

# --- Column-width formatting moves from column B to column A ---
# (column A now holds "name", which used to be column B's custom width)
$ws.Columns.Item(1).ColumnWidth = 7.125
$ws.Columns.Item(2).ColumnWidth = 9.71

# --- Selection moves to C10 ---
$ws.Range("C10").Select()
